# This workbook is a scraped "cryptos" price table refresh (GitHub Actions bot
# commit). The table is rebuilt on every run: the coin that disappeared from
# the top of the list (WazirX) gets bumped down, a brand-new coin ("One")
# gets inserted above it, every row below shifts down by one, and most of the
# "Price" (column D) values are refreshed with the latest quote. A handful of
# other, unrelated rows further down the sheet also get their price touched
# up, and a couple of "Volume(1h)" (column E) labels change too.
#
# All of the numeric-looking prices are stored as *text* (inlineStr) in the
# original file, not as real numbers - trailing zeros like "0.05900" or
# "0.0002000" are meaningful and must survive. If we just do
#   $ws.Range("D5").Value = "0.05900"
# Excel will happily "smart type" that into the number 0.059 and we lose the
# formatting/trailing zeros (and the cell representation changes type).
# Prefixing the value with a leading apostrophe forces Excel to keep it as
# literal text, and then resetting the cell style back to "Normal" afterwards
# keeps the cell's style index identical to the untouched cells around it
# (the apostrophe/text coercion otherwise tags the cell with a different
# number-format style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$Address,
        [string]$Text
    )
    $range = $ws.Range($Address)
    # Leading apostrophe => store literally as text, never auto-convert to a number/date.
    $range.Value = "'" + $Text
    # Restore the default style so we don't leave a stray "text" number format
    # on the cell (the source cells have no explicit style at all).
    $range.Style = "Normal"
}

# --- Rows 3-9: price refresh only -----------------------------------------
Set-TextCell "D3" "23.02"
Set-TextCell "D4" "5.420"
Set-TextCell "D5" "0.05900"
Set-TextCell "D6" "3.442"
Set-TextCell "D7" "6.518"
Set-TextCell "D8" "0.8103"
Set-TextCell "D9" "0.9579"

# --- Rows 10-18: "One" is newly inserted at rank 9, pushing WazirX and
# everything through CoinExToken down by one row. -------------------------
Set-TextCell "B10" "One"
Set-TextCell "C10" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell "D10" "0.01129"
Set-TextCell "E10" "9OneONEBestin24h"

Set-TextCell "B11" "WazirX"
Set-TextCell "C11" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell "D11" "0.1420"
Set-TextCell "E11" "10WazirXWRX"

Set-TextCell "B12" "MandalaExchangeToken"
Set-TextCell "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D12" "0.07425"
Set-TextCell "E12" "11MandalaExchangeTokenMDX"

Set-TextCell "B13" "LiechtensteinCryptoassetsExchange"
Set-TextCell "C13" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell "D13" "0.03264"
Set-TextCell "E13" "12LiechtensteinCryptoassetsExchangeLCX"

Set-TextCell "B14" "BitrueCoin"
Set-TextCell "C14" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell "D14" "0.03052"
Set-TextCell "E14" "13BitrueCoinBTR"

Set-TextCell "B15" "BitMartToken"
Set-TextCell "C15" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell "D15" "0.09340"
Set-TextCell "E15" "14BitMartTokenBMX"

Set-TextCell "B16" "MCDex"
Set-TextCell "C16" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell "D16" "3.846"
Set-TextCell "E16" "15MCDexMCB"

Set-TextCell "B17" "BitForexToken"
Set-TextCell "C17" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell "D17" "0.001592"
Set-TextCell "E17" "16BitForexTokenBF"

Set-TextCell "B18" "CoinExToken"
Set-TextCell "C18" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell "D18" "0.04669"
Set-TextCell "E18" "17CoinExTokenCET"

# --- Remaining rows: price refresh only (coin/link/rank label unchanged) --
Set-TextCell "D19" "0.005881"
Set-TextCell "D20" "0.001265"
Set-TextCell "D21" "0.004902"
Set-TextCell "D22" "0.00006801"
Set-TextCell "D23" "3.574"
Set-TextCell "D24" "2.131"
Set-TextCell "D27" "0.0002284"
Set-TextCell "D40" "0.03936"
Set-TextCell "D41" "0.006185"
Set-TextCell "D43" "0.003000"
Set-TextCell "D44" "0.009911"
Set-TextCell "E44" "43LocalTradersLCT"
Set-TextCell "D45" "0.00005200"
Set-TextCell "D46" "0.00000000750"
Set-TextCell "D47" "0.7201"
Set-TextCell "D48" "0.002392"
Set-TextCell "D49" "0.00002100"
Set-TextCell "D50" "0.0002000"
